$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 31 de Agosto de 2020 a las 00:45"

# Update country case-count rows (data refresh + rank-swap renames)
# Row 4
$ws.Cells.Item(4, 2).Value = 6170330
$ws.Cells.Item(4, 3).Value = 30960
$ws.Cells.Item(4, 4).Value = 3422602
$ws.Cells.Item(4, 5).Value = 2560527
$ws.Cells.Item(4, 7).Value = 346
$ws.Cells.Item(4, 8).Value = 187201

# Row 5
$ws.Cells.Item(5, 2).Value = 3862311
$ws.Cells.Item(5, 3).Value = 15346
$ws.Cells.Item(5, 4).Value = 3031559
$ws.Cells.Item(5, 5).Value = 709924
$ws.Cells.Item(5, 7).Value = 330
$ws.Cells.Item(5, 8).Value = 120828

# Row 10
$ws.Cells.Item(10, 2).Value = 607938
$ws.Cells.Item(10, 3).Value = 8024
$ws.Cells.Item(10, 4).Value = 450621
$ws.Cells.Item(10, 5).Value = 137953
$ws.Cells.Item(10, 7).Value = 300
$ws.Cells.Item(10, 8).Value = 19364

# Row 14
$ws.Cells.Item(14, 2).Value = 408426
$ws.Cells.Item(14, 3).Value = 7187
$ws.Cells.Item(14, 5).Value = 105962
$ws.Cells.Item(14, 7).Value = 104
$ws.Cells.Item(14, 8).Value = 8457

# Row 34
$ws.Cells.Item(34, 2).Value = 98727
$ws.Cells.Item(34, 3).Value = 230
$ws.Cells.Item(34, 4).Value = 72120
$ws.Cells.Item(34, 5).Value = 21208
$ws.Cells.Item(34, 7).Value = 23
$ws.Cells.Item(34, 8).Value = 5399

# Row 36
$ws.Cells.Item(36, 2).Value = 92065
$ws.Cells.Item(36, 3).Value = 728
$ws.Cells.Item(36, 4).Value = 65747
$ws.Cells.Item(36, 5).Value = 24323
$ws.Cells.Item(36, 7).Value = 12
$ws.Cells.Item(36, 8).Value = 1995

# Row 43
$ws.Cells.Item(43, 2).Value = 73912
$ws.Cells.Item(43, 3).Value = 233
$ws.Cells.Item(43, 4).Value = 61681
$ws.Cells.Item(43, 5).Value = 9491
$ws.Cells.Item(43, 7).Value = 12
$ws.Cells.Item(43, 8).Value = 2740

# Row 47
$ws.Cells.Item(47, 1).Value = "Japon"
$ws.Cells.Item(47, 2).Value = 67264
$ws.Cells.Item(47, 3).Value = 841
$ws.Cells.Item(47, 4).Value = 56164
$ws.Cells.Item(47, 5).Value = 9836
$ws.Cells.Item(47, 7).Value = 9
$ws.Cells.Item(47, 8).Value = 1264

# Row 48
$ws.Cells.Item(48, 1).Value = "Polonia"
$ws.Cells.Item(48, 2).Value = 66870
$ws.Cells.Item(48, 3).Value = 631
$ws.Cells.Item(48, 4).Value = 46192
$ws.Cells.Item(48, 5).Value = 18645
$ws.Cells.Item(48, 7).Value = 1
$ws.Cells.Item(48, 8).Value = 2033

# Row 53
$ws.Cells.Item(53, 2).Value = 53865
$ws.Cells.Item(53, 3).Value = 138
$ws.Cells.Item(53, 4).Value = 41513
$ws.Cells.Item(53, 5).Value = 11339
$ws.Cells.Item(53, 7).Value = 2
$ws.Cells.Item(53, 8).Value = 1013

# Row 58
$ws.Cells.Item(58, 1).Value = "Argelia"
$ws.Cells.Item(58, 2).Value = 44146
$ws.Cells.Item(58, 3).Value = 365
$ws.Cells.Item(58, 4).Value = 30978
$ws.Cells.Item(58, 5).Value = 11667
$ws.Cells.Item(58, 7).Value = 10
$ws.Cells.Item(58, 8).Value = 1501

# Row 59
$ws.Cells.Item(59, 1).Value = "Kirguistan"
$ws.Cells.Item(59, 2).Value = 43820
$ws.Cells.Item(59, 3).Value = 108
$ws.Cells.Item(59, 4).Value = 38198
$ws.Cells.Item(59, 5).Value = 4564
$ws.Cells.Item(59, 7).Value = 1
$ws.Cells.Item(59, 8).Value = 1058

# Row 83
$ws.Cells.Item(83, 2).Value = 16190
$ws.Cells.Item(83, 3).Value = 26
$ws.Cells.Item(83, 4).Value = 11313
$ws.Cells.Item(83, 5).Value = 4264
$ws.Cells.Item(83, 7).Value = 8
$ws.Cells.Item(83, 8).Value = 613

# Row 90
$ws.Cells.Item(90, 2).Value = 10643
$ws.Cells.Item(90, 3).Value = 32
$ws.Cells.Item(90, 5).Value = 1031

# Row 95
$ws.Cells.Item(95, 1).Value = "Guinea"
$ws.Cells.Item(95, 2).Value = 9371
$ws.Cells.Item(95, 3).Value = 120
$ws.Cells.Item(95, 4).Value = 8387
$ws.Cells.Item(95, 5).Value = 925
$ws.Cells.Item(95, 7).Value = 0
$ws.Cells.Item(95, 8).Value = 59

# Row 96
$ws.Cells.Item(96, 1).Value = "Malasia"
$ws.Cells.Item(96, 2).Value = 9334
$ws.Cells.Item(96, 3).Value = 17
$ws.Cells.Item(96, 4).Value = 9048
$ws.Cells.Item(96, 5).Value = 160
$ws.Cells.Item(96, 7).Value = 1
$ws.Cells.Item(96, 8).Value = 126

# Row 121
$ws.Cells.Item(121, 1).Value = "Cabo Verde"
$ws.Cells.Item(121, 2).Value = 3852
$ws.Cells.Item(121, 3).Value = 74
$ws.Cells.Item(121, 4).Value = 2867
$ws.Cells.Item(121, 5).Value = 945
$ws.Cells.Item(121, 8).Value = 40

# Row 122
$ws.Cells.Item(122, 1).Value = "Mozambique"
$ws.Cells.Item(122, 2).Value = 3821
$ws.Cells.Item(122, 3).Value = 61
$ws.Cells.Item(122, 4).Value = 2100
$ws.Cells.Item(122, 5).Value = 1698
$ws.Cells.Item(122, 7).Value = 1
$ws.Cells.Item(122, 8).Value = 23

# Row 123
$ws.Cells.Item(123, 2).Value = 3685
$ws.Cells.Item(123, 3).Value = 113
$ws.Cells.Item(123, 4).Value = 1562
$ws.Cells.Item(123, 5).Value = 2047

# Row 158
$ws.Cells.Item(158, 2).Value = 1234
$ws.Cells.Item(158, 3).Value = 50
$ws.Cells.Item(158, 4).Value = 687
$ws.Cells.Item(158, 5).Value = 511
$ws.Cells.Item(158, 7).Value = 1
$ws.Cells.Item(158, 8).Value = 36
